# Generate Report for Archive
# 1) Status text: "Ready for handoff" -> "In Translation" on every sheet
#    that carries the localization Status column.
# 2) Narrow the Status-related columns (Overview!E:F, zh-cn!C, de-de!C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Replace the status text wherever it appears ---
foreach ($cellRef in @("E2","F2","E3","F3","E4","F4")) {
    if ($overview.Range($cellRef).Value2 -eq "Ready for handoff") {
        $overview.Range($cellRef).Value = "In Translation"
    }
}

foreach ($cellRef in @("C2","C3","C4")) {
    if ($zhcn.Range($cellRef).Value2 -eq "Ready for handoff") {
        $zhcn.Range($cellRef).Value = "In Translation"
    }
    if ($dede.Range($cellRef).Value2 -eq "Ready for handoff") {
        $dede.Range($cellRef).Value = "In Translation"
    }
}

# --- Narrow the columns (17.2159881591797 -> 13.4101848602295 raw width) ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
